$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 21 (shifts the "Summa Psicológica" row and
# everything below it down by one) to make room for the new journal entry
# ("Peace and Conflict: Journal of Peace Psychology") in the multicols
# journal list that lives in column E.
$ws.Rows(21).Insert()

# Fill in the new row with the review entry, matching the style already
# used by the surrounding list rows (left/top aligned, wrapped text).
$ws.Range("E21").Value = '\href{https://www.apa.org/pubs/journals/pac}{Peace and Conflict: Journal of \newline Peace Psychology}'
$ws.Range("E21").HorizontalAlignment = -4131
$ws.Range("E21").VerticalAlignment = -4160
$ws.Range("E21").WrapText = $true

# Reflect the scrolled viewport (the author had scrolled the sheet so
# column C / row 6 is the top-left visible cell) without touching the
# current selection.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 6
